# Add a new "2022" column (P) to the table, mirroring the style of the
# existing 2021 column (O), then fill in the new year's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header (row 4): new year label, formatted like the other year headers.
$ws.Range("O4").Copy()
$ws.Range("P4").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("P4").Value = 2022

# Data rows: copy number formatting from column O, then set the new values.
$values = @{
    5  = 11.4
    6  = 12.6
    7  = 9.8
    8  = 11.4
    9  = 5.4
    10 = 4.7
    11 = 3.4
    12 = 17.7
    13 = 20.5
    14 = 8.4
    16 = 12.9
    17 = 10.5
}

foreach ($row in $values.Keys) {
    $src = $ws.Range("O$row")
    $dst = $ws.Range("P$row")
    $src.Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats
    $dst.Value = $values[$row]
}

# Move the active selection to Q4, matching the saved view state.
$ws.Range("Q4").Select()
